$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "UniformA-HW45.xpc" to "UniformA"
$ws.Name = "UniformA"

# Add a new row (16) following the same pattern as row 15, using the
# "HexGrid-60degTilt5degRes" label (already present as a shared string)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16:P16").Value = 1

# Column A carries the bold/centered/bordered style used throughout the sheet
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
